$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(40,2).Value = 754
$ws.Cells.Item(41,2).Value = 811
$ws.Cells.Item(42,2).Value = 541
$ws.Range("C41:C42").FillDown()
$ws.Range("D41:D42").FillDown()
Write-Host "done"
